$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend header row with new month columns ---
$ws.Range("GZ1").Value = "2024-10"
$ws.Range("HA1").Value = "2024-11"
# Copy style (bold/border/center) from the preceding header cell GY1 so GZ1:HA1 match formatting
$ws.Range("GY1").Copy()
$ws.Range("GZ1:HA1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Fill in / extend data values for rows 2-33 ---
# Row 2
$ws.Range("GY2").Value = 42.507
$ws.Range("GZ2").Value = 42.417
$ws.Range("HA2").Value = 41.048
# Row 3
$ws.Range("GY3").Value = 1663.697
$ws.Range("GZ3").Value = 1660.154
$ws.Range("HA3").Value = 1606.601
# Row 4
$ws.Range("GY4").Value = 0.944
$ws.Range("GZ4").Value = 0.967
$ws.Range("HA4").Value = 0.9360000000000001
# Row 5
$ws.Range("GY5").Value = 36.946
$ws.Range("GZ5").Value = 37.853
$ws.Range("HA5").Value = 36.632
# Row 6
$ws.Range("GY6").Value = 1293.64
$ws.Range("GZ6").Value = 1023.216
$ws.Range("HA6").Value = 925.816
# Row 7
$ws.Range("GY7").Value = 50632.225
$ws.Range("GZ7").Value = 40047.99
$ws.Range("HA7").Value = 36235.842
# Row 8
$ws.Range("GY8").Value = 888.443
$ws.Range("GZ8").Value = 378.323
$ws.Range("HA8").Value = 529.657
# Row 9
$ws.Range("GY9").Value = 34773.036
$ws.Range("GZ9").Value = 14807.353
$ws.Range("HA9").Value = 20730.408
# Row 10
$ws.Range("GY10").Value = 0
$ws.Range("GZ10").Value = 0
$ws.Range("HA10").Value = 0
# Row 11
$ws.Range("GY11").Value = 0
$ws.Range("GZ11").Value = 0
$ws.Range("HA11").Value = 0
# Row 12
$ws.Range("GY12").Value = 67.44499999999999
$ws.Range("GZ12").Value = 153.243
$ws.Range("HA12").Value = -499.813
# Row 13
$ws.Range("GY13").Value = 2639.763
$ws.Range("GZ13").Value = 5997.825
$ws.Range("HA13").Value = -19562.333
# Row 14
$ws.Range("GY14").Value = 0
$ws.Range("GZ14").Value = 0
$ws.Range("HA14").Value = 0
# Row 15
$ws.Range("GY15").Value = 0
$ws.Range("GZ15").Value = 0
$ws.Range("HA15").Value = 0
# Row 16
$ws.Range("GY16").Value = 376.43
$ws.Range("GZ16").Value = 534.51
$ws.Range("HA16").Value = 0
# Row 17
$ws.Range("GY17").Value = 14733.239
$ws.Range("GZ17").Value = 20920.369
$ws.Range("HA17").Value = 0
# Row 18
$ws.Range("GY18").Value = 381.203
$ws.Range("GZ18").Value = 535.034
$ws.Range("HA18").Value = 937.956
# Row 19
$ws.Range("GY19").Value = 14920.069
$ws.Range("GZ19").Value = 20940.819
$ws.Range("HA19").Value = 36711
# Row 20
$ws.Range("GY20").Value = 23.026
$ws.Range("GZ20").Value = 83.524
$ws.Range("HA20").Value = 0
# Row 21
$ws.Range("GY21").Value = 901.236
$ws.Range("GZ21").Value = 3269.088
$ws.Range("HA21").Value = 0
# Row 22
$ws.Range("GY22").Value = 16.474
$ws.Range("GZ22").Value = 10.345
$ws.Range("HA22").Value = 0
# Row 23
$ws.Range("GY23").Value = 644.775
$ws.Range("GZ23").Value = 404.909
$ws.Range("HA23").Value = 0
# Row 24
$ws.Range("GY24").Value = 0
$ws.Range("GZ24").Value = 0
$ws.Range("HA24").Value = 0
# Row 25
$ws.Range("GY25").Value = 0
$ws.Range("GZ25").Value = 0
$ws.Range("HA25").Value = 0
# Row 26
$ws.Range("GY26").Value = 0
$ws.Range("GZ26").Value = 0
$ws.Range("HA26").Value = 0
# Row 27
$ws.Range("GY27").Value = 0
$ws.Range("GZ27").Value = 0
$ws.Range("HA27").Value = 0
# Row 28
# GZ28 remains blank (empty string cell) - no action needed
$ws.Range("HA28").Value = 0
# Row 29
# GZ29 remains blank (empty string cell) - no action needed
$ws.Range("HA29").Value = 0
# Row 30
$ws.Range("GY30").Value = 71.072
$ws.Range("GZ30").Value = 140.356
$ws.Range("HA30").Value = 0
# Row 31
$ws.Range("GY31").Value = 2781.703
$ws.Range("GZ31").Value = 5493.426
$ws.Range("HA31").Value = 0
# Row 32
$ws.Range("GY32").Value = 4.773
$ws.Range("GZ32").Value = 0.524
$ws.Range("HA32").Value = 937.956
# Row 33
$ws.Range("GY33").Value = 186.83
$ws.Range("GZ33").Value = 20.45
$ws.Range("HA33").Value = 36711

Write-Output "edit complete"
